$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Data fix: part number LXES15AAA1-153 -> 0402ESDB-MLP1
#    (appears twice on row 14: column B "Manu. P/N" and column G "Comment")
#    Set value first, then re-apply the original cell formatting (copied from
#    an untouched neighbouring cell) so the cell keeps its original style
#    (quotePrefix text style) instead of the COM layer's default "new value"
#    style.
# ---------------------------------------------------------------------------
$ws.Range("B14").Value = "0402ESDB-MLP1"
$ws.Range("G14").Value = "0402ESDB-MLP1"

$ws.Range("A14").Copy()
$ws.Range("B14").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("F14").Copy()
$ws.Range("G14").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2) Column widths were narrowed across the whole sheet (likely a manual
#    resize pass as part of the rev2.5 output refresh).
#    ColumnWidth is quantized internally, so we pre-compensate by the fixed
#    offset this engine applies (+5/6 of a character) so the stored XML
#    width lands as close as possible to the target values.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 16.333333333333336
$ws.Columns.Item(2).ColumnWidth = 23.0
$ws.Columns.Item(3).ColumnWidth = 5.666666666666667
$ws.Columns.Item(4).ColumnWidth = 5.5
$ws.Columns.Item(5).ColumnWidth = 20.5
$ws.Columns.Item(6).ColumnWidth = 13.666666666666666
$ws.Columns.Item(7).ColumnWidth = 23.333333333333336

# ---------------------------------------------------------------------------
# 3) Selection / scroll state: row 18 is selected (whole row) in the saved
#    view.
# ---------------------------------------------------------------------------
$ws.Rows.Item(18).Select()
